$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Export as TSV"

$ws8 = $wb.Worksheets.Item(8)
$ws8.Name = "transposition_...se_source list"

# --- 2. Freeze top row on the main sheet ---
$ws1.Activate()
$ws1.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- 3. Update data validations: add errorTitle/error, and fix the Y column formula ---

function Set-DV($range, $title, $msg) {
    $r = $ws1.Range($range)
    $r.Validation.ErrorTitle = $title
    $r.Validation.ErrorMessage = $msg
}

Set-DV "I2:I1048576"  "Value must come from list" "Value must be one of: sequence."
Set-DV "J2:J1048576"  "Value must come from list" "Value must be one of: SNARE-seq2 / scATACseq / sciATACseq / snATACseq."
Set-DV "K2:K1048576"  "Value must come from list" "Value must be one of: DNA."
Set-DV "L2:L1048576"  "Not a boolean" 'The values in this column must be "TRUE" or "FALSE".'
Set-DV "O2:O1048576"  "Not a boolean" 'The values in this column must be "TRUE" or "FALSE".'
Set-DV "R2:R1048576"  "Value must come from list" "Value must be one of: whole cell / nucleus / cell-cell multimer / spatially encoded cell barcoding."
Set-DV "T2:T1048576"  "Value must come from list" "Value must be one of: none / FACS."
Set-DV "V2:V1048576"  "Not an integer" "The values in this column must be integers."
Set-DV "W2:W1048576"  "Not a number" "The values in this column must be numbers."
Set-DV "X2:X1048576"  "Value must come from list" "Value must be one of: SNARE-Seq2-AC / scATACseq / bulkATACseq / snATACseq / sciATACseq."

$rngY = $ws1.Range("Y2:Y1048576")
$rngY.Validation.Formula1 = "='transposition_...se_source list'!`$A`$1:`$A`$3"
$rngY.Validation.ErrorTitle = "Value must come from list"
$rngY.Validation.ErrorMessage = "Value must be one of: 10X snATAC / In-house / Nextera."

Set-DV "AB2:AB1048576" "Value must come from list" "Value must be one of: single-end / paired-end."
Set-DV "AG2:AG1048576" "Not an integer" "The values in this column must be integers."
Set-DV "AH2:AH1048576" "Not an integer" "The values in this column must be integers."
Set-DV "AI2:AI1048576" "Not a number" "The values in this column must be numbers."
Set-DV "AJ2:AJ1048576" "Value must come from list" "Value must be one of: ng."
Set-DV "AK2:AK1048576" "Not a number" "The values in this column must be numbers."
Set-DV "AN2:AN1048576" "Not a number" "The values in this column must be numbers."
Set-DV "AO2:AO1048576" "Not a number" "The values in this column must be numbers."

Write-Host "Done"
